$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 29
$ws.Range("H29").Value = 1252
$ws.Range("J29").Value = 2004
$ws.Range("L29").Value = 6012
$ws.Range("N29").Value = -6574
# Row 38
$ws.Range("H38").Value = 5421.6
$ws.Range("I38").Value = 869.3333
$ws.Range("J38").Value = 12250
$ws.Range("K38").Value = 2607.9999
$ws.Range("L38").Value = 36750
$ws.Range("M38").Value = -2235.9999
$ws.Range("N38").Value = -37494
# Row 43
$ws.Range("H43").Value = 499.66666
$ws.Range("I43").Value = 499.66666
$ws.Range("K43").Value = 499.66666
$ws.Range("M43").Value = -430.66666
# Row 58
$ws.Range("H58").Value = 2957.75
# Row 133
$ws.Range("H133").Value = 79554.5
$ws.Range("J133").Value = 79554.5
$ws.Range("L133").Value = 79554.5
$ws.Range("N133").Value = -89674.5
# Row 137
$ws.Range("H137").Value = 498572.28
$ws.Range("I137").Value = 1032756.7
$ws.Range("J137").Value = 12950.091
$ws.Range("K137").Value = 3098270.1
$ws.Range("L137").Value = 38850.273
$ws.Range("M137").Value = -3095720.1
$ws.Range("N137").Value = -43950.273

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2539.4443
$ws.Range("I32").Value = 2707.0344
$ws.Range("K32").Value = 2707.0344
$ws.Range("M32").Value = -2420.0344
# Row 61
$ws.Range("H61").Value = 120000
$ws.Range("I61").Value = 120000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 120000
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -119788
$ws.Range("N61").ClearContents()
# Row 122
$ws.Range("H122").Value = 641199.1
$ws.Range("I122").Value = 4394.923
$ws.Range("K122").Value = 13184.769
$ws.Range("M122").Value = -10734.769
# Row 132
$ws.Range("H132").Value = 3242.1
$ws.Range("I132").Value = 2453.2144
$ws.Range("K132").Value = 7359.6432
$ws.Range("M132").Value = -4829.6432
# Row 136
$ws.Range("H136").Value = 120000
$ws.Range("I136").Value = 120000
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 360000
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -357450
$ws.Range("N136").ClearContents()

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1927.0264
$ws.Range("I20").Value = 883.6070999999999
$ws.Range("J20").Value = 4848.6
$ws.Range("K20").Value = 883.6070999999999
$ws.Range("L20").Value = 4848.6
$ws.Range("M20").Value = -636.6070999999999
$ws.Range("N20").Value = -5342.6
# Row 22
$ws.Range("H22").Value = 709.625
$ws.Range("I22").Value = 709.625
$ws.Range("K22").Value = 709.625
$ws.Range("M22").Value = -536.625
# Row 140
$ws.Range("H140").Value = 86999.5
$ws.Range("J140").Value = 86999.5
$ws.Range("L140").Value = 86999.5
$ws.Range("N140").Value = -97359.5
# Row 141
$ws.Range("H141").Value = 117499.75
$ws.Range("J141").Value = 117499.75
$ws.Range("L141").Value = 117499.75
$ws.Range("N141").Value = -127859.75

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 17
$ws.Range("H17").Value = 95000
$ws.Range("I17").Value = 95000
$ws.Range("K17").Value = 95000
$ws.Range("M17").Value = -94826
# Row 31
$ws.Range("H31").Value = 3387.5557
$ws.Range("I31").Value = 2118
$ws.Range("J31").Value = 3592.3225
$ws.Range("K31").Value = 2118
$ws.Range("L31").Value = 3592.3225
$ws.Range("M31").Value = -1823
$ws.Range("N31").Value = -4182.3225
# Row 34
$ws.Range("H34").Value = 3387.5557
$ws.Range("I34").Value = 2118
$ws.Range("J34").Value = 3592.3225
$ws.Range("K34").Value = 2118
$ws.Range("L34").Value = 3592.3225
$ws.Range("M34").Value = -1916
$ws.Range("N34").Value = -3996.3225
# Row 99
$ws.Range("H99").Value = 1002542.4
$ws.Range("I99").Value = 1667570.6
$ws.Range("K99").Value = 1667570.6
$ws.Range("M99").Value = -1666072.6
# Row 122
$ws.Range("H122").Value = 2220.3333
$ws.Range("I122").Value = 2219.5
$ws.Range("K122").Value = 6658.5
$ws.Range("M122").Value = -4208.5
# Row 126
$ws.Range("H126").Value = 1002542.4
$ws.Range("I126").Value = 1667570.6
$ws.Range("K126").Value = 5002711.800000001
$ws.Range("M126").Value = -5000241.800000001
# Row 132
$ws.Range("H132").Value = 53859.3
$ws.Range("I132").Value = 26357.8
$ws.Range("J132").Value = 81360.8
$ws.Range("K132").Value = 79073.39999999999
$ws.Range("L132").Value = 244082.4
$ws.Range("M132").Value = -76543.39999999999
$ws.Range("N132").Value = -249142.4

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 92
$ws.Range("H92").Value = 2050
$ws.Range("J92").Value = 2050
$ws.Range("L92").Value = 6150
$ws.Range("N92").Value = -8646

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 9267.764999999999
$ws.Range("I70").Value = 9562.75
$ws.Range("J70").Value = 8559.799999999999
$ws.Range("K70").Value = 9562.75
$ws.Range("L70").Value = 8559.799999999999
$ws.Range("M70").Value = -9292.75
$ws.Range("N70").Value = -9099.799999999999
# Row 73
$ws.Range("H73").Value = 9267.764999999999
$ws.Range("I73").Value = 9562.75
$ws.Range("J73").Value = 8559.799999999999
$ws.Range("K73").Value = 9562.75
$ws.Range("L73").Value = 8559.799999999999
$ws.Range("M73").Value = -8626.75
$ws.Range("N73").Value = -10431.8
# Row 113
$ws.Range("H113").Value = 2999.5
$ws.Range("I113").Value = 2999.5
$ws.Range("K113").Value = 2999.5
$ws.Range("M113").Value = -829.5
# Row 126
$ws.Range("H126").Value = 26623.8
$ws.Range("I126").Value = 53999.668
$ws.Range("K126").Value = 161999.004
$ws.Range("M126").Value = -159529.004
# Row 132
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 23338
$ws.Range("I22").Value = 29047.5
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 29047.5
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -28752.5
$ws.Range("N22").Value = -1090
# Row 27
$ws.Range("H27").Value = 23338
$ws.Range("I27").Value = 29047.5
$ws.Range("J27").Value = 500
$ws.Range("K27").Value = 29047.5
$ws.Range("L27").Value = 500
$ws.Range("M27").Value = -28940.5
$ws.Range("N27").Value = -714
# Row 40
$ws.Range("H40").Value = 98335.5
$ws.Range("I40").Value = 234833.33
$ws.Range("J40").Value = 16436.8
$ws.Range("K40").Value = 234833.33
$ws.Range("L40").Value = 16436.8
$ws.Range("M40").Value = -234697.33
$ws.Range("N40").Value = -16708.8
# Row 61
$ws.Range("H61").Value = 3515
$ws.Range("I61").Value = 3551.5
$ws.Range("J61").Value = 3466.3333
$ws.Range("K61").Value = 3551.5
$ws.Range("L61").Value = 3466.3333
$ws.Range("M61").Value = -3349.5
$ws.Range("N61").Value = -3870.3333
# Row 93
$ws.Range("H93").Value = 4107.0625
$ws.Range("I93").Value = 4107.0625
$ws.Range("K93").Value = 4107.0625
$ws.Range("M93").Value = -2859.0625
# Row 100
$ws.Range("H100").Value = 8054.75
$ws.Range("I100").Value = 9510.777
$ws.Range("J100").Value = 3686.6667
$ws.Range("K100").Value = 9510.777
$ws.Range("L100").Value = 3686.6667
$ws.Range("M100").Value = -8969.777
$ws.Range("N100").Value = -4768.6667
# Row 113
$ws.Range("H113").Value = 3515
$ws.Range("I113").Value = 3551.5
$ws.Range("J113").Value = 3466.3333
$ws.Range("K113").Value = 3551.5
$ws.Range("L113").Value = 3466.3333
$ws.Range("M113").Value = -1381.5
$ws.Range("N113").Value = -7806.3333

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 100638.305
$ws.Range("J62").Value = 4458.0835
$ws.Range("L62").Value = 4458.0835
$ws.Range("N62").Value = -5706.0835
# Row 65
$ws.Range("H65").Value = 100638.305
$ws.Range("J65").Value = 4458.0835
$ws.Range("L65").Value = 22290.4175
$ws.Range("N65").Value = -28530.4175
# Row 122
$ws.Range("H122").Value = 4442.5713
$ws.Range("I122").Value = 3745.0908
$ws.Range("J122").Value = 7000
$ws.Range("K122").Value = 11235.2724
$ws.Range("L122").Value = 21000
$ws.Range("M122").Value = -8785.2724
$ws.Range("N122").Value = -25900
# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").ClearContents()
$ws.Range("N123").ClearContents()
# Row 126
$ws.Range("H126").Value = 20013.217
$ws.Range("I126").Value = 22390.2
$ws.Range("K126").Value = 67170.60000000001
$ws.Range("M126").Value = -64700.60000000001
